$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# Columns E (codeforiati:category-name) and F (codeforiati:group-code) need to
# swap places, for the header row and for every data row. Column F typically
# holds digit-only codes (e.g. "110"), so after the swap those land in column
# E; force a Text number format while writing so Excel's auto type-detection
# doesn't turn them into numbers, then restore the default "Normal" style so
# no stray cell formatting is left behind.
for ($r = 1; $r -le $lastRow; $r++) {
    $eVal = $ws.Cells.Item($r, 5).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2

    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value2 = $fVal
    $ws.Cells.Item($r, 5).Style = "Normal"

    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value2 = $eVal
    $ws.Cells.Item($r, 6).Style = "Normal"
}
